$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 10.48 = 42468.66 pesos`n✅ 42468.66 pesos = 10.44 = 955.26 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 95.40000000000001
$ws2.Range("O10").Value = 4051.51

$ws2.Range("N12").Value = 4068
$ws2.Range("O12").Value = 91.503
